$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 178343
$ws.Range("C4").Value = 168300
$ws.Range("C7").Value = 5.63
$ws.Range("C8").Value = 64.93000000000001
